$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 22:53:22"
$wsZhCn.Range("H2").Value = "2016-03-18 22:53:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 22:53:25"
$wsDeDe.Range("H2").Value = "2016-03-18 22:53:48"
